$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the F-column values (rows 2-5) to 10
$ws.Range("F2:F5").Value = 10

# Move the selection/active cell to M1 (also drops the scrolled topLeftCell)
$ws.Range("M1").Select()
